$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1
$ws.Range("B1").Value = "yYow"

$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("E3").Value = "asdasd"

$ws.Range("B2").Value = "asd"
$ws.Range("C2").Value = "asd"
$ws.Range("D2").Value = "asda"
$ws.Range("E2").Value = "sd"

$ws.Range("C3").Value = "Dasdas"
$ws.Range("D3").Value = "asdas"

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "agasfags"
$ws.Range("C4").Value = "dfefgxg"
$ws.Range("D4").Value = "sfdsghs"
$ws.Range("E4").Value = "sfsdf"

$ws.Range("A5").Value = 5
$ws.Range("A6").Value = 6
$ws.Range("A7").Value = 7
$ws.Range("A8").Value = 8
$ws.Range("A9").Value = 9
$ws.Range("A10").Value = 10

$ws.Range("A11").Select() | Out-Null
